# Products_Lazada_filtered.xlsx - "Final Result of Lazada"
#
# Each worksheet has columns:
#   A Name | B Price (text, e.g. "$14.90") | C From | D Sold (text, e.g. "43 sold" / "2.2K sold")
#   E (was "sold_cleaned")  | F (was "price_cleaned")
#
# The edit re-derives the cleaned numeric helper columns so that:
#   E = price_cleaned  (precise numeric price parsed from column B)
#   F = sold_cleaned   (numeric sold count parsed from column D)
# and swaps the column headers to match ("price_cleaned" in E1, "sold_cleaned" in F1).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {

        $priceText = [string]$ws.Cells.Item($r, 2).Value2
        $soldText  = [string]$ws.Cells.Item($r, 4).Value2

        # price_cleaned: strip currency symbol / thousands separators -> number
        $priceClean = [double]($priceText.Trim() -replace '[\$,]', '')

        # sold_cleaned: strip the " sold" suffix, expand "K" (thousands) suffix -> number
        $soldTrim = ($soldText -replace '\s*sold\s*$', '').Trim()
        if ($soldTrim -match '^(?<num>[0-9.]+)K$') {
            $soldClean = [double]$Matches['num'] * 1000
        } else {
            $soldClean = [double]$soldTrim
        }

        $ws.Cells.Item($r, 5).Value2 = $priceClean
        $ws.Cells.Item($r, 6).Value2 = $soldClean
    }

    $ws.Cells.Item(1, 5).Value2 = "price_cleaned"
    $ws.Cells.Item(1, 6).Value2 = "sold_cleaned"
}
